# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Update column G ("K") values for rows 2-37 on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 4
    3  = 2
    4  = 0
    5  = 6
    6  = 9
    7  = 4
    8  = 4
    9  = 3
    10 = 3
    11 = 5
    12 = 3
    13 = 3
    14 = 2
    15 = 4
    16 = 0
    17 = 1
    18 = 4
    19 = 6
    20 = 4
    21 = 1
    22 = 8
    23 = 5
    24 = 6
    25 = 6
    26 = 3
    27 = 4
    28 = 6
    29 = 3
    30 = 6
    31 = 4
    32 = 7
    33 = 6
    34 = 3
    35 = 3
    36 = 3
    37 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
